$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B ("Polarity"), shifting Polarity/Review right.
$ws.Range("B:B").Insert()

# New column B header + values mirror column A ("Unnamed: 0" -> "Unnamed: 0.1")
$ws.Range("A1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$ws.Range("B1").Value = "Unnamed: 0.1"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 2).Value = $ws.Cells.Item($r, 1).Value2
}

# Lowercase the review text, now in column D.
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $cell.Value = $cell.Value2.ToString().ToLower()
}
